# aggiornamento fino a 6 gennaio 2022
# Append new daily-data rows (465-491) to the single data sheet, mirroring
# the existing layout: col A = date serial (styled like the rows above),
# col B = nuovi pos., col C = somma mobile 7gg., col D = somma mobile 7gg.
# per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(465, 44539, 0, 1, 48.07692307692308),
  @(466, 44540, 1, 2, 96.15384615384616),
  @(467, 44541, 0, 2, 96.15384615384616),
  @(468, 44542, 1, 3, 144.2307692307692),
  @(469, 44543, 0, 3, 144.2307692307692),
  @(470, 44544, 0, 2, 96.15384615384616),
  @(471, 44545, 0, 2, 96.15384615384616),
  @(472, 44546, 0, 2, 96.15384615384616),
  @(473, 44547, 3, 4, 192.3076923076923),
  @(474, 44548, 2, 6, 288.4615384615385),
  @(475, 44550, 1, 6, 288.4615384615385),
  @(476, 44551, 0, 6, 288.4615384615385),
  @(477, 44552, 0, 6, 288.4615384615385),
  @(478, 44553, 0, 6, 288.4615384615385),
  @(479, 44554, 0, 6, 288.4615384615385),
  @(480, 44555, 1, 4, 192.3076923076923),
  @(481, 44556, 0, 2, 96.15384615384616),
  @(482, 44557, 0, 1, 48.07692307692308),
  @(483, 44558, 0, 1, 48.07692307692308),
  @(484, 44559, 5, 6, 288.4615384615385),
  @(485, 44560, 5, 11, 528.8461538461539),
  @(486, 44561, 2, 13, 625),
  @(487, 44562, 3, 15, 721.1538461538462),
  @(488, 44563, 2, 17, 817.3076923076924),
  @(489, 44564, 2, 19, 913.4615384615386),
  @(490, 44565, 0, 19, 913.4615384615386),
  @(491, 44566, 2, 16, 769.2307692307693)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Column A carries the same date/time number format + styling as the rest
# of the column (style index used by A2:A464); replicate it onto the newly
# added cells by copying the format from the last pre-existing row.
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Output "Appended rows 465-491 (through 2022-01-05)."
